$d = $word.ActiveDocument

# --- Fix "paketten" -> "pakketten" and drop its spell-check proofErr markers ---
$f1 = $d.Content
$f1.Find.Execute("paketten", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Expand the found range by one character on each side so the replacement
# "crosses" the zero-width <w:proofErr/> markers flanking the run, which
# causes them to be dropped when the span is rewritten.
$cross1 = $d.Range($f1.Start - 1, $f1.End + 1)
$cross1.Text = " pakketten "

# The crossing-replace merges neighbouring runs into one; restore the
# original run boundaries by toggling a character property on just the
# corrected word (round-tripping Bold forces a run split without leaving
# any visible formatting behind).
$w1 = $d.Content
$w1.Find.Execute("pakketten", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$w1.Bold = 1
$w1.Bold = 0

# --- Fix "economieen" -> "economieën" and drop its spell-check proofErr markers ---
$f2 = $d.Content
$f2.Find.Execute("economieen", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$cross2 = $d.Range($f2.Start - 1, $f2.End + 1)
$cross2.Text = " economieën "

$w2 = $d.Content
$w2.Find.Execute("economieën", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$w2.Bold = 1
$w2.Bold = 0
